$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link swaps (rows 16-24), plain text cells ---
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('B18').Value = 'BitpandaEcosystemToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('B19').Value = 'MCDex'
$ws.Range('C19').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('B20').Value = 'ProBitToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('B21').Value = 'ZBToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('B22').Value = 'CoinExToken'
$ws.Range('C22').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('B23').Value = 'BitKan'
$ws.Range('C23').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('B24').Value = 'HotbitToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'

# --- Price / Volume(1h) cells: force Text format so numeric-looking strings
#     (e.g. "309.46", "0.36%") are not auto-converted to numbers/percentages ---
$priceVolumeCells = @(
    'D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5',
    'D6', 'E6', 'D7', 'E7', 'D9', 'E9', 'D10', 'E10',
    'D11', 'E11', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14',
    'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'D18', 'E18',
    'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22',
    'D23', 'E23', 'D24', 'E24', 'D25', 'E25', 'D38', 'E38',
    'D39', 'E39', 'D40', 'E40', 'D41', 'E41', 'D42', 'E42',
    'E43', 'D44', 'E44', 'D45', 'E45', 'D46', 'E46', 'E47',
    'D48', 'E48', 'D50', 'E50', 'D51', 'E51'
)
foreach ($addr in $priceVolumeCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '309.46'
$ws.Range('E2').Value = '0.36%'
$ws.Range('D3').Value = '41.16'
$ws.Range('E3').Value = '0.75%'
$ws.Range('D4').Value = '5.217'
$ws.Range('E4').Value = '2.13%'
$ws.Range('D5').Value = '0.07685'
$ws.Range('E5').Value = '0.82%'
$ws.Range('D6').Value = '1.644'
$ws.Range('E6').Value = '2.57%'
$ws.Range('D7').Value = '0.9147'
$ws.Range('E7').Value = '1.46%'
$ws.Range('D9').Value = '0.1244'
$ws.Range('E9').Value = '10.51%'
$ws.Range('D10').Value = '0.1823'
$ws.Range('E10').Value = '1.63%'
$ws.Range('D11').Value = '0.09175'
$ws.Range('E11').Value = '0.29%'
$ws.Range('D12').Value = '0.04225'
$ws.Range('E12').Value = '1.42%'
$ws.Range('D13').Value = '0.1053'
$ws.Range('E13').Value = '0.09%'
$ws.Range('D14').Value = '0.001249'
$ws.Range('E14').Value = '-0.09%'
$ws.Range('D15').Value = '0.005890'
$ws.Range('E15').Value = '2.47%'
$ws.Range('D16').Value = '3.350'
$ws.Range('E16').Value = '0.03%'
$ws.Range('D17').Value = '4.316'
$ws.Range('E17').Value = '1.41%'
$ws.Range('D18').Value = '0.3336'
$ws.Range('E18').Value = '0.65%'
$ws.Range('D19').Value = '7.379'
$ws.Range('E19').Value = '10.88%'
$ws.Range('D20').Value = '0.1402'
$ws.Range('E20').Value = '2.68%'
$ws.Range('D21').Value = '0.2820'
$ws.Range('E21').Value = '0.64%'
$ws.Range('D22').Value = '0.04022'
$ws.Range('E22').Value = '-1.35%'
$ws.Range('D23').Value = '0.001265'
$ws.Range('E23').Value = '1.56%'
$ws.Range('D24').Value = '0.004097'
$ws.Range('E24').Value = '-0.11%'
$ws.Range('D25').Value = '0.0001301'
$ws.Range('E25').Value = '0.00%'
$ws.Range('D38').Value = '0.02552'
$ws.Range('E38').Value = '6.60%'
$ws.Range('D39').Value = '0.05348'
$ws.Range('E39').Value = '3.10%'
$ws.Range('D40').Value = '0.007828'
$ws.Range('E40').Value = '0.53%'
$ws.Range('D41').Value = '0.1314'
$ws.Range('E41').Value = '1.18%'
$ws.Range('D42').Value = '0.006678'
$ws.Range('E42').Value = '-5.40%'
$ws.Range('E43').Value = '-4.62%'
$ws.Range('D44').Value = '0.008055'
$ws.Range('E44').Value = '4.13%'
$ws.Range('D45').Value = '0.3070'
$ws.Range('E45').Value = '-0.25%'
$ws.Range('D46').Value = '0.00006723'
$ws.Range('E46').Value = '-3.55%'
$ws.Range('E47').Value = '0.00%'
$ws.Range('D48').Value = '0.2215'
$ws.Range('E48').Value = '374.04%'
$ws.Range('D50').Value = '0.00002102'
$ws.Range('E50').Value = '0.00%'
$ws.Range('D51').Value = '0.0002002'
$ws.Range('E51').Value = '0.00%'
